# Rename transcript speaker labels in the "Speaker" column (D):
#   "Davis"   -> "T"
#   "Student" -> "S"
# All other speaker names are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D = Speaker
    $val = $cell.Text

    if ($val -eq "Davis") {
        $cell.Value = "T"
    }
    elseif ($val -eq "Student") {
        $cell.Value = "S"
    }
}
